# Generate Report for Handoff
# Update the "latest generate/handoff" timestamps for the file
# 953e0f08-aecb-4b53-a766-72bc9d9918db.md (row 7 on every sheet) to reflect
# a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 7
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-09-05 02:47:01"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for row 7
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-09-05 02:46:56"

# de-de sheet: "Latest Handoff Datetime" column (H) for row 7
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-09-05 02:47:01"
